$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Selig 1223
$ws.Range("A4").Value = "Selig 1223"
$ws.Range("B4").Value = -6.7
$ws.Range("C4").Value = 0.121
$ws.Range("D4").Value = 0.49
$ws.Range("E4").Formula = "=(1.2135-1.1864)/0.25"
$ws.Range("F4").Value = -0.29
$ws.Range("G4").Value = 1.1864
$ws.Range("H4").Value = 0.023

# Row 5: NACA 0012
$ws.Range("A5").Value = "NACA 0012"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0.12
$ws.Range("D5").Value = 0
$ws.Range("E5").Formula = "=(1.6-1)/(2+4)"
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0

# Match the final selection state left behind in the source workbook
$ws.Range("D8").Select()
